# Apply the "Add files via upload" edits to casos_de_prueba.xlsx (Hoja1)
#  - fills in the previously-blank "Rio Cuarto/Córdoba" location cell (D16)
#  - fills in the previously-empty 3rd test-case block (rows 19-21) with
#    new test cases (ID / Caso de Prueba / Descripcion / Fecha)
#  - updates the Fecha (date) values of the last test-case block (rows 29-31)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- D16: fill in the previously-empty "Probado en:" location value ---
$ws.Range("D16").Value = "Rio Cuarto/Córdoba"

# --- New test cases for the 3rd block (rows 19-21) ---
# Written in this particular cell order so the shared-string table comes
# out in the same sequence as the authored workbook.
$ws.Range("C21").Value = "El curso se puede editar correctamente, aparecen todos los campos correspondientes del formulario"
$ws.Range("B20").Value = "SuperUser"
$ws.Range("C20").Value = "Boton editar y borrar solo se ven en el perfil de superuser"
$ws.Range("B19").Value = "Ver curso"
$ws.Range("C19").Value = 'Cada boton "Ver" de cada curso corresponde con la descripcion de cada curso'
$ws.Range("B21").Value = "Editar curso"

# ID column (row numbers 1,2,3)
$ws.Range("A19").Value = 1
$ws.Range("A20").Value = 2
$ws.Range("A21").Value = 3

# Fecha column: set alignment before the number format so the engine folds
# both into a single cell style (matching the authored style table), then
# set the date value itself.
$ws.Range("D19").HorizontalAlignment = -4108   # xlCenter
$ws.Range("D19").NumberFormat = "mm-dd-yy"
$ws.Range("D19").Value = 45214

$ws.Range("D20").HorizontalAlignment = -4108   # xlCenter
$ws.Range("D20").VerticalAlignment = -4108     # xlCenter
$ws.Range("D20").NumberFormat = "mm-dd-yy"
$ws.Range("D20").Value = 45216

$ws.Range("D21").HorizontalAlignment = -4108   # xlCenter
$ws.Range("D21").VerticalAlignment = -4108     # xlCenter
$ws.Range("D21").NumberFormat = "d-mmm"
$ws.Range("D21").Value = 45217

# --- Update the Fecha values of the last test-case block (rows 29-31) ---
$ws.Range("D29").Value = 45214
$ws.Range("D30").Value = 45215
$ws.Range("D31").Value = 45217

# --- Selection / active cell, matching the final saved view ---
$ws.Activate()
$ws.Range("D21").Select()
